$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Tabela FIPE" row with a test entry (bad URL) to test the alarm
$ws.Range("B7").Value = "https://brasilapi.com.br/api/feriados/v1/{year}"
$ws.Range("A7").Value = "TESTE ALARME"

# Update the active selection to A8, matching the saved view state
$ws.Range("A8").Select()
